$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.609.73"
$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "2.375.30"
$ws.Range("E3").Value = "  -4.22%  "
$ws.Range("D5").Value = "'311.52"
$ws.Range("E5").Value = "  -2.46%  "
$ws.Range("D6").Value = "'86.91"
$ws.Range("E6").Value = "  -6.33%  "
$ws.Range("D7").Value = "'0.533"
$ws.Range("E7").Value = "  -3.46%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -4.19%  "
$ws.Range("D10").Value = "'0.0833"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").Value = "'30.36"
$ws.Range("E11").Value = "  -8.63%  "
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Value = "2.742.58"
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").Value = "'6.51"
$ws.Range("E14").Value = "  -5.72%  "
$ws.Range("D15").Value = "'14.95"
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("D16").Value = "2.366.44"
$ws.Range("E16").Value = "  -4.45%  "
$ws.Range("D17").Value = "'0.759"
$ws.Range("E17").Value = "  -4.54%  "
$ws.Range("D18").Value = "40.578.38"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").Value = "0.0₃0911"
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").Value = "'6.14"
$ws.Range("E20").Value = "  -5.02%  "
$ws.Range("D21").Value = "'68.37"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("D22").Value = "'10.81"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").Value = "'235.16"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").Value = "'2.59"
$ws.Range("E24").Value = "  -6.01%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -7.93%  "
$ws.Range("D27").Value = "'23.71"
$ws.Range("E27").Value = "  -5.48%  "
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "'9.24"
$ws.Range("E29").Value = "  -5.07%  "
$ws.Range("D30").Value = "'34.32"
$ws.Range("E30").Value = "  -7.23%  "
$ws.Range("D31").Value = "'153.35"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'5.22"
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("D34").Value = "'0.0730"
$ws.Range("E34").Value = "  -4.85%  "
$ws.Range("D35").Value = "'2.42"
$ws.Range("E35").Value = "  -5.31%  "
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("D37").Value = "'2.80"
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("D38").Value = "'15.94"
$ws.Range("E38").Value = "  -7.50%  "
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("E40").Value = "  -7.84%  "
$ws.Range("D43").Value = "1.969.19"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  -5.58%  "
$ws.Range("D45").Value = "'17.64"
$ws.Range("E45").Value = "  -6.80%  "
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("E47").Value = "  -9.54%  "
$ws.Range("D48").Value = "2.601.44"
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("D49").Value = "'93.29"
$ws.Range("E49").Value = "  -5.19%  "
$ws.Range("D50").Value = "'71.54"
$ws.Range("E50").Value = "  -5.60%  "
$ws.Range("D51").Value = "'50.40"
$ws.Range("E51").Value = "  -3.91%  "

# Row 41 and 42 swap order/content (RenderToken now ranks above ApeXProtocol)
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'3.82"
$ws.Range("E41").Value = "  -5.15%  "

$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "'2.39"
$ws.Range("E42").Value = "  -3.80%  "
